# Apply the "Upload new version with timestamp" changes to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (first item line: سرنجات 3 سم) value updates
$ws.Range("H7").Value = "0:0"
$ws.Range("Q7").Value = "11:0"

# P7 keeps its numeric-looking display ("22.0000") but is stored as literal
# text (matches the source file, which already stored "2.0000" as text even
# though the cell carries a 0.00 numeric style). Temporarily force a text
# format so COM doesn't coerce the value to a real number, then restore the
# original 0.00 number format.
$p7Format = $ws.Range("P7").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "22.0000"
$ws.Range("P7").NumberFormat = $p7Format

# Row 9 (totals row) numeric update
$ws.Range("N9").Value = 25

# Row 10 footer timestamp update (minute rolled from :44 to :45)
$ws.Range("A10").Value = "Wednesday, 17 September, 2025 10:45 PM"
